# Sheet1 is the active/selected sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "selenium_GF1" value in D2 with the new test data name.
$ws.Range("D2").Value = "Name_AlreadyExist_Search"

# Move the active selection to C2, matching the refreshed sheet view.
$ws.Range("C2").Select()
